$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2, "face/face010.jpg", "lügen", "face"),
  @(3, "house/house023.jpg", "lernen", "house"),
  @(4, "house/house008.jpg", "danken", "house"),
  @(5, "face/face020.jpg", "heißen", "face"),
  @(6, "house/house005.jpg", "sparen", "house"),
  @(7, "house/house003.jpg", "achten", "house"),
  @(8, "face/face013.jpg", "zögern", "face"),
  @(9, "face/face004.jpg", "rechnen", "face"),
  @(10, "face/face007.jpg", "süßen", "face"),
  @(11, "face/face005.jpg", "bergen", "face"),
  @(12, "house/house013.jpg", "wachsen", "house"),
  @(13, "face/face023.jpg", "spüren", "face"),
  @(14, "house/house027.jpg", "betteln", "house"),
  @(15, "house/house004.jpg", "leeren", "house"),
  @(16, "house/house016.jpg", "lächeln", "house"),
  @(17, "face/face003.jpg", "holen", "face"),
  @(18, "house/house026.jpg", "hören", "house"),
  @(19, "face/face014.jpg", "dienen", "face"),
  @(20, "face/face019.jpg", "hacken", "face"),
  @(21, "house/house028.jpg", "ärgern", "house"),
  @(22, "face/face017.jpg", "wählen", "face"),
  @(23, "face/face030.jpg", "öffnen", "face"),
  @(24, "house/house002.jpg", "prüfen", "house"),
  @(25, "face/face002.jpg", "lassen", "face"),
  @(26, "house/house030.jpg", "frischen", "house"),
  @(27, "face/face001.jpg", "planen", "face"),
  @(28, "house/house010.jpg", "stoßen", "house"),
  @(29, "face/face021.jpg", "zielen", "face"),
  @(30, "face/face027.jpg", "duschen", "face"),
  @(31, "house/house029.jpg", "angeln", "house"),
  @(32, "house/house011.jpg", "hassen", "house"),
  @(33, "house/house031.jpg", "proben", "house")
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
}
